$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The sheet is protected; unprotect it so the cell values can be updated,
# then restore protection afterwards.
$ws.Unprotect()

# Update the confidential/model-holdings-as-of date string in A10 (2021-07-13 -> 2021-07-14)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-14 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) columns for rows 2-7
$ws.Range("D2").Value = 0.2644018982318107
$ws.Range("E2").Value = -0.002596212583995117

$ws.Range("D3").Value = 0.530496757377827
$ws.Range("E3").Value = 0.0005278437582474371

$ws.Range("D4").Value = 0.05216073145215178
$ws.Range("E4").Value = -0.007614678899082739

$ws.Range("D5").Value = 0.09575153745086774
$ws.Range("E5").Value = 0.0008321775312067903

$ws.Range("D6").Value = 0.05718907548734273
$ws.Range("E6").Value = -0.009840515778757886

$ws.Range("D7").Value = 0.9999999999999999
$ws.Range("E7").Value = -0.001286699076079367

# Restore sheet protection (the sheet was protected before this edit, with
# objects and scenarios also protected).
$ws.Protect($null, $true, $true, $true)
